$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently lists 9 "Bereta" bookings (room 1) across rows 2-10.
# The new layout keeps the first "Bereta" booking (rows 2-5, which already
# repeats the 23rd-25th) and replaces rows 6-10 with a brand-new block of
# "Stanuszek" bookings (room 2) that reuses the same dates/times as rows
# 2-6. The four "Bereta" bookings that used to occupy rows 7-10 simply move
# down to rows 11-14 (which were blank before) so they keep existing.

# 1) Move the trailing "Bereta" bookings (old rows 7-10) down to rows 11-14.
$ws.Range("A7:E10").Copy($ws.Range("A11:E14"))

# 2) Re-seed rows 6-10 with a copy of the date/time pattern from rows 2-6.
$ws.Range("A2:E6").Copy($ws.Range("A6:E10"))

# 3) Turn that new block into the "Stanuszek" / room 2 bookings.
$ws.Range("D6:D10").Value = "Stanuszek"
$ws.Range("E6:E10").Value = 2

# 4) The copy above dragged along number formats from the source rows;
#    fix the two dates whose format needs to differ from their source.
$ws.Range("A6").NumberFormat = "dd/mm/yyyy"
$ws.Range("A10").NumberFormat = "yyyy-mm-dd"

# 5) Normalize the "przewodniczacy" cells in the touched rows back to the
#    default (unstyled) look used elsewhere in the sheet.
$ws.Range("D6:D14").Style = "Normal"
